$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired state for rows 2-17 (columns A: Player, B: Position, C: Team)
$data = @(
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Jalen Johnson", "PF", "Atlanta Hawks"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
